# "added restauration and demolition methods"
#
# Tweak four numeric parameter values on the single data sheet of
# parameter_scenarios.xlsx (soe / sme scenario table):
#
#   - soe / restoration_rate : 2020 value            0.5 -> 0.4
#   - soe / restoration_rate : divergence column      "0.7, 0.8, 1, 0.4" -> 0
#   - soe / restoration_ab   : divergence column      0.3 -> 0
#   - soe / demolition_rate_min : 2020 value          0.005 -> 0.1

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: soe / restoration_rate
$ws.Range("C2").Value = 0.4
$ws.Range("I2").Value = 0

# Row 7: soe / restoration_ab
$ws.Range("I7").Value = 0

# Row 8: soe / demolition_rate_min
$ws.Range("C8").Value = 0.1
